# Add a header row to the Colleges sheet ("Added in headers for reference").
# Inserts a brand-new row 1 (shifting all existing data down by one row) and
# fills it with column headers: ID, Name, Degree, Site, Region.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 1 -- everything else (values,
# shared strings, per-cell styles) shifts down automatically.
[void]$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Degree"
$ws.Range("D1").Value = "Site"
$ws.Range("E1").Value = "Region"

# Match the post-edit selection (first data row, column A).
[void]$ws.Range("A2").Select()
